$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" field text from 2020/12/21 to
#    2020/12/22 everywhere it is cached: the slide master's Date
#    Placeholder and every slide layout's Date Placeholder.
# ---------------------------------------------------------------------
$oldDate = "2020/12/21"
$newDate = "2020/12/22"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1 shape edits.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# -- Shape id=6 ("矩形 5"): nudge its position and append the project
#    number "PRJ2020-002" after the existing title text.
$rectShape = $slide.Shapes.Item(7)

$rectShape.Left = -35790 / 12700
$rectShape.Top = 28900755 / 12700

$titleRange = $rectShape.TextFrame.TextRange
$lenBefore = $titleRange.Length

$null = $titleRange.InsertAfter(" ")
$lenAfterSpace = $titleRange.Length

$null = $titleRange.InsertAfter("PRJ2020-002")
$lenAfterPrj = $titleRange.Length

$prjRun = $titleRange.Characters($lenAfterSpace + 1, $lenAfterPrj - $lenAfterSpace)
$prjRun.Font.Name = "Times New Roman"

# -- Shape id=43 ("文字方塊 42" / "Conclusions"): nudge its position.
$conclusionsShape = $slide.Shapes.Item(14)
$conclusionsShape.Left = 14613107 / 12700
$conclusionsShape.Top = 19992844 / 12700
